$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.138.70"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.900.98"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'324.56"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.4626"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "'0.3903"
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("D9").Value = "'0.07875"
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").Value = "'0.9919"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").Value = "'21.96"
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").Value = "1.858.49"
$ws.Range("E12").Value = "  -2.91%  "
$ws.Range("D13").Value = "'5.787"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").Value = "'7.065"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").Value = "'0.06988"
$ws.Range("D16").Value = "'87.94"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "'0.000009941"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("D21").Value = "29.131.22"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "'5.318"
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").Value = "'11.10"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "2.142.70"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "'2.104"
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("D26").Value = "'155.78"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "'19.39"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").Value = "'5.914"
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("D29").Value = "'118.59"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").Value = "'1.881"
$ws.Range("E30").Value = "  -5.63%  "
$ws.Range("D31").Value = "'0.09318"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").Value = "'0.9005"
$ws.Range("D33").Value = "'5.248"
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("D34").Value = "'1.327"
$ws.Range("E34").Value = "  -1.73%  "
$ws.Range("D35").Value = "'3.161"
$ws.Range("E35").Value = "  -2.95%  "
$ws.Range("D36").Value = "'0.05808"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("D38").Value = "'0.02084"
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("D39").Value = "'0.9991"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").Value = "'7.722"
$ws.Range("E40").Value = "  -2.92%  "
$ws.Range("D41").Value = "'0.5688"
$ws.Range("E41").Value = "  -0.94%  "
$ws.Range("D42").Value = "'0.1796"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "'9.751"
$ws.Range("E43").Value = "  -2.00%  "
$ws.Range("D44").Value = "'2.239"
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("D46").Value = "'0.5358"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("D48").Value = "'1.848"
$ws.Range("E48").Value = "  -1.41%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").Value = "'112.92"
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("E51").Value = "  -0.33%  "
